$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Nlgn2"
$ws.Range("C2").Value = "Nrxn1"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.580781333333333
$ws.Range("H2").Value = 7.742344
$ws.Range("I2").Value = 0.09250274102763278
$ws.Range("J2").Value = 0.0925027410276328
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2191816666666667
$ws.Range("N2").Value = 0.657545
$ws.Range("O2").Value = 0.3167322406056181
$ws.Range("P2").Value = 0.3167322406056181
$ws.Range("Q2").Value = 0.5656599539422223
$ws.Range("R2").Value = 5.09093958548
$ws.Range("S2").Value = 0.02929860042784336
$ws.Range("T2").Value = 0.02929860042784337

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Nlgn2"
$ws.Range("C3").Value = "Nrxn1"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.580781333333333
$ws.Range("H3").Value = 7.742344
$ws.Range("I3").Value = 0.09250274102763278
$ws.Range("J3").Value = 0.0925027410276328
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.4728276666666666
$ws.Range("N3").Value = 1.418483
$ws.Range("O3").Value = 0.6832677593943819
$ws.Range("P3").Value = 0.6832677593943819
$ws.Range("Q3").Value = 1.220264816016889
$ws.Range("R3").Value = 10.982383344152
$ws.Range("S3").Value = 0.06320414059978942
$ws.Range("T3").Value = 0.06320414059978942

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Nlgn2"
$ws.Range("C4").Value = "Nrxn1"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 20.80139666666667
$ws.Range("H4").Value = 62.40418999999999
$ws.Range("I4").Value = 0.7455828140172008
$ws.Range("J4").Value = 0.7455828140172009
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.2191816666666667
$ws.Range("N4").Value = 0.657545
$ws.Range("O4").Value = 0.3167322406056181
$ws.Range("P4").Value = 0.3167322406056181
$ws.Range("Q4").Value = 4.559284790394445
$ws.Range("R4").Value = 41.03356311355
$ws.Range("S4").Value = 0.2361501152407098
$ws.Range("T4").Value = 0.2361501152407099

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Nlgn2"
$ws.Range("C5").Value = "Nrxn1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 20.80139666666667
$ws.Range("H5").Value = 62.40418999999999
$ws.Range("I5").Value = 0.7455828140172008
$ws.Range("J5").Value = 0.7455828140172009
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.4728276666666666
$ws.Range("N5").Value = 1.418483
$ws.Range("O5").Value = 0.6832677593943819
$ws.Range("P5").Value = 0.6832677593943819
$ws.Range("Q5").Value = 9.835475849307777
$ws.Range("R5").Value = 88.51928264376998
$ws.Range("S5").Value = 0.509432698776491
$ws.Range("T5").Value = 0.509432698776491

$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Nlgn2"
$ws.Range("C6").Value = "Nrxn1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.305159
$ws.Range("H6").Value = 0.915477
$ws.Range("I6").Value = 0.01093778987962227
$ws.Range("J6").Value = 0.01093778987962227
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2191816666666667
$ws.Range("N6").Value = 0.657545
$ws.Range("O6").Value = 0.3167322406056181
$ws.Range("P6").Value = 0.3167322406056181
$ws.Range("Q6").Value = 0.06688525821833334
$ws.Range("R6").Value = 0.6019673239650001
$ws.Range("S6").Value = 0.003464350695846214
$ws.Range("T6").Value = 0.003464350695846215

$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Nlgn2"
$ws.Range("C7").Value = "Nrxn1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.305159
$ws.Range("H7").Value = 0.915477
$ws.Range("I7").Value = 0.01093778987962227
$ws.Range("J7").Value = 0.01093778987962227
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.4728276666666666
$ws.Range("N7").Value = 1.418483
$ws.Range("O7").Value = 0.6832677593943819
$ws.Range("P7").Value = 0.6832677593943819
$ws.Range("Q7").Value = 0.1442876179323333
$ws.Range("R7").Value = 1.298588561391
$ws.Range("S7").Value = 0.007473439183776053
$ws.Range("T7").Value = 0.007473439183776054

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Nlgn2"
$ws.Range("C8").Value = "Nrxn1"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 4.212175
$ws.Range("H8").Value = 12.636525
$ws.Range("I8").Value = 0.150976655075544
$ws.Range("J8").Value = 0.150976655075544
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.2191816666666667
$ws.Range("N8").Value = 0.657545
$ws.Range("O8").Value = 0.3167322406056181
$ws.Range("P8").Value = 0.3167322406056181
$ws.Range("Q8").Value = 0.9232315367916668
$ws.Range("R8").Value = 8.309083831125001
$ws.Range("S8").Value = 0.0478191742412186
$ws.Range("T8").Value = 0.04781917424121861

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Nlgn2"
$ws.Range("C9").Value = "Nrxn1"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 4.212175
$ws.Range("H9").Value = 12.636525
$ws.Range("I9").Value = 0.150976655075544
$ws.Range("J9").Value = 0.150976655075544
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.4728276666666666
$ws.Range("N9").Value = 1.418483
$ws.Range("O9").Value = 0.6832677593943819
$ws.Range("P9").Value = 0.6832677593943819
$ws.Range("Q9").Value = 1.991632876841667
$ws.Range("R9").Value = 17.924695891575
$ws.Range("S9").Value = 0.1031574808343254
$ws.Range("T9").Value = 0.1031574808343254
